$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.193268
$ws.Range("H2").Value = 0.579804
$ws.Range("I2").Value = 0.1207017725010034
$ws.Range("J2").Value = 0.1207017725010034
$ws.Range("M2").Value = 1.684857333333333
$ws.Range("N2").Value = 5.054572
$ws.Range("O2").Value = 0.03007292173576635
$ws.Range("P2").Value = 0.03007292173576635
$ws.Range("Q2").Value = 0.3256290070986667
$ws.Range("R2").Value = 2.930661063888
$ws.Range("S2").Value = 0.00362985495779095
$ws.Range("T2").Value = 0.003629854957790949
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.193268
$ws.Range("H3").Value = 0.579804
$ws.Range("I3").Value = 0.1207017725010034
$ws.Range("J3").Value = 0.1207017725010034
$ws.Range("O3").Value = 0.8213562764116752
$ws.Range("P3").Value = 0.8213562764116752
$ws.Range("Q3").Value = 8.893629661666667
$ws.Range("R3").Value = 80.042666955
$ws.Range("S3").Value = 0.09913915841771329
$ws.Range("T3").Value = 0.09913915841771327
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.193268
$ws.Range("H4").Value = 0.579804
$ws.Range("I4").Value = 0.1207017725010034
$ws.Range("J4").Value = 0.1207017725010034
$ws.Range("M4").Value = 8.323787333333334
$ws.Range("N4").Value = 24.971362
$ws.Range("O4").Value = 0.1485708018525584
$ws.Range("P4").Value = 0.1485708018525584
$ws.Range("Q4").Value = 1.608721730338667
$ws.Range("R4").Value = 14.478495573048
$ws.Range("S4").Value = 0.01793275912549916
$ws.Range("T4").Value = 0.01793275912549916
$ws.Range("I5").Value = 0.7607038292883184
$ws.Range("J5").Value = 0.7607038292883183
$ws.Range("M5").Value = 1.684857333333333
$ws.Range("N5").Value = 5.054572
$ws.Range("O5").Value = 0.03007292173576635
$ws.Range("P5").Value = 0.03007292173576635
$ws.Range("Q5").Value = 2.052225311150667
$ws.Range("R5").Value = 18.470027800356
$ws.Range("S5").Value = 0.02287658672228536
$ws.Range("T5").Value = 0.02287658672228536
$ws.Range("I6").Value = 0.7607038292883184
$ws.Range("J6").Value = 0.7607038292883183
$ws.Range("O6").Value = 0.8213562764116752
$ws.Range("P6").Value = 0.8213562764116752
$ws.Range("S6").Value = 0.6248088646763558
$ws.Range("T6").Value = 0.6248088646763557
$ws.Range("I7").Value = 0.7607038292883184
$ws.Range("J7").Value = 0.7607038292883183
$ws.Range("M7").Value = 8.323787333333334
$ws.Range("N7").Value = 24.971362
$ws.Range("O7").Value = 0.1485708018525584
$ws.Range("P7").Value = 0.1485708018525584
$ws.Range("Q7").Value = 10.13871424728067
$ws.Range("R7").Value = 91.248428225526
$ws.Range("S7").Value = 0.1130183778896772
$ws.Range("T7").Value = 0.1130183778896771
$ws.Range("G8").Value = 0.1898936666666667
$ws.Range("H8").Value = 0.569681
$ws.Range("I8").Value = 0.1185943982106783
$ws.Range("J8").Value = 0.1185943982106783
$ws.Range("M8").Value = 1.684857333333333
$ws.Range("N8").Value = 5.054572
$ws.Range("O8").Value = 0.03007292173576635
$ws.Range("P8").Value = 0.03007292173576635
$ws.Range("Q8").Value = 0.3199437368368889
$ws.Range("R8").Value = 2.879493631532
$ws.Range("S8").Value = 0.003566480055690037
$ws.Range("T8").Value = 0.003566480055690037
$ws.Range("G9").Value = 0.1898936666666667
$ws.Range("H9").Value = 0.569681
$ws.Range("I9").Value = 0.1185943982106783
$ws.Range("J9").Value = 0.1185943982106783
$ws.Range("O9").Value = 0.8213562764116752
$ws.Range("P9").Value = 0.8213562764116752
$ws.Range("Q9").Value = 8.738352683472224
$ws.Range("R9").Value = 78.64517415125
$ws.Range("S9").Value = 0.09740825331760616
$ws.Range("T9").Value = 0.09740825331760616
$ws.Range("G10").Value = 0.1898936666666667
$ws.Range("H10").Value = 0.569681
$ws.Range("I10").Value = 0.1185943982106783
$ws.Range("J10").Value = 0.1185943982106783
$ws.Range("M10").Value = 8.323787333333334
$ws.Range("N10").Value = 24.971362
$ws.Range("O10").Value = 0.1485708018525584
$ws.Range("P10").Value = 0.1485708018525584
$ws.Range("Q10").Value = 1.580634497280222
$ws.Range("R10").Value = 14.225710475522
$ws.Range("S10").Value = 0.01761966483738209
$ws.Range("T10").Value = 0.01761966483738209
